$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.293291666666667
$ws.Range("H2").Value = 9.879875
$ws.Range("I2").Value = 0.161585160668464
$ws.Range("J2").Value = 0.1670404328173815
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.927965666666666
$ws.Range("N2").Value = 5.783897
$ws.Range("O2").Value = 0.4937378937520984
$ws.Range("P2").Value = 0.4937378937520984
$ws.Range("Q2").Value = 6.349353263652778
$ws.Range("R2").Value = 57.14417937287499
$ws.Range("S2").Value = 0.07978071689004185
$ws.Range("T2").Value = 0.08247419147069285

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.293291666666667
$ws.Range("H3").Value = 9.879875
$ws.Range("I3").Value = 0.161585160668464
$ws.Range("J3").Value = 0.1670404328173815
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.976870666666667
$ws.Range("N3").Value = 5.930612
$ws.Range("O3").Value = 0.5062621062479017
$ws.Range("P3").Value = 0.5062621062479017
$ws.Range("Q3").Value = 6.510411692611111
$ws.Range("R3").Value = 58.5937052335
$ws.Range("S3").Value = 0.08180444377842221
$ws.Range("T3").Value = 0.08456624134668869

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.839091333333334
$ws.Range("H4").Value = 8.517274
$ws.Range("I4").Value = 0.1392998482012507
$ws.Range("J4").Value = 0.1440027465311282
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.927965666666666
$ws.Range("N4").Value = 5.783897
$ws.Range("O4").Value = 0.4937378937520984
$ws.Range("P4").Value = 0.4937378937520984
$ws.Range("Q4").Value = 5.473670615197555
$ws.Range("R4").Value = 49.263035536778
$ws.Range("S4").Value = 0.06877761365087254
$ws.Range("T4").Value = 0.07109961276679655

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.839091333333334
$ws.Range("H5").Value = 8.517274
$ws.Range("I5").Value = 0.1392998482012507
$ws.Range("J5").Value = 0.1440027465311282
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.976870666666667
$ws.Range("N5").Value = 5.930612
$ws.Range("O5").Value = 0.5062621062479017
$ws.Range("P5").Value = 0.5062621062479017
$ws.Range("Q5").Value = 5.612516376854223
$ws.Range("R5").Value = 50.51264739168801
$ws.Range("S5").Value = 0.07052223455037815
$ws.Range("T5").Value = 0.0729031337643317

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 6.486612666666666
$ws.Range("H6").Value = 19.459838
$ws.Range("I6").Value = 0.3182652664949994
$ws.Range("J6").Value = 0.3290102113717155
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.927965666666666
$ws.Range("N6").Value = 5.783897
$ws.Range("O6").Value = 0.4937378937520984
$ws.Range("P6").Value = 0.4937378937520984
$ws.Range("Q6").Value = 12.50596651429844
$ws.Range("R6").Value = 112.553698628686
$ws.Range("S6").Value = 0.1571396223336913
$ws.Range("T6").Value = 0.1624448087856035

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 6.486612666666666
$ws.Range("H7").Value = 19.459838
$ws.Range("I7").Value = 0.3182652664949994
$ws.Range("J7").Value = 0.3290102113717155
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.976870666666667
$ws.Range("N7").Value = 5.930612
$ws.Range("O7").Value = 0.5062621062479017
$ws.Range("P7").Value = 0.5062621062479017
$ws.Range("Q7").Value = 12.82319430676178
$ws.Range("R7").Value = 115.408748760856
$ws.Range("S7").Value = 0.1611256441613081
$ws.Range("T7").Value = 0.166565402586112

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 5.765308999999999
$ws.Range("H8").Value = 17.295927
$ws.Range("I8").Value = 0.2828745447897899
$ws.Range("J8").Value = 0.292424664487945
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.927965666666666
$ws.Range("N8").Value = 5.783897
$ws.Range("O8").Value = 0.4937378937520984
$ws.Range("P8").Value = 0.4937378937520984
$ws.Range("Q8").Value = 11.11531780972433
$ws.Range("R8").Value = 100.037860287519
$ws.Range("S8").Value = 0.1396658819405945
$ws.Range("T8").Value = 0.144381137925442

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 5.765308999999999
$ws.Range("H9").Value = 17.295927
$ws.Range("I9").Value = 0.2828745447897899
$ws.Range("J9").Value = 0.292424664487945
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.976870666666667
$ws.Range("N9").Value = 5.930612
$ws.Range("O9").Value = 0.5062621062479017
$ws.Range("P9").Value = 0.5062621062479017
$ws.Range("Q9").Value = 11.39727024636933
$ws.Range("R9").Value = 102.575432217324
$ws.Range("S9").Value = 0.1432086628491954
$ws.Range("T9").Value = 0.148043526562503

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("G10").Value = 1.996847
$ws.Range("H10").Value = 3.993694
$ws.Range("I10").Value = 0.09797517984549615
$ws.Range("J10").Value = 0.06752194479182982
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.927965666666666
$ws.Range("N10").Value = 5.783897
$ws.Range("O10").Value = 0.4937378937520984
$ws.Range("P10").Value = 0.4937378937520984
$ws.Range("Q10").Value = 3.849852457586333
$ws.Range("R10").Value = 23.099114745518
$ws.Range("S10").Value = 0.04837405893689831
$ws.Range("T10").Value = 0.03333814280356352

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("G11").Value = 1.996847
$ws.Range("H11").Value = 3.993694
$ws.Range("I11").Value = 0.09797517984549615
$ws.Range("J11").Value = 0.06752194479182982
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.976870666666667
$ws.Range("N11").Value = 5.930612
$ws.Range("O11").Value = 0.5062621062479017
$ws.Range("P11").Value = 0.5062621062479017
$ws.Range("Q11").Value = 3.947508260121333
$ws.Range("R11").Value = 23.685049560728
$ws.Range("S11").Value = 0.04960112090859785
$ws.Range("T11").Value = 0.0729031337643317
